$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 - this shifts every existing data row
# (2..41) down by one (3..42), so the old last row (2025-11-21) simply
# ends up at row 42 untouched, and the whole table's dimension grows to
# A1:D42 automatically.
$ws.Rows.Item(2).Insert()

# The inserted row inherits formatting from the row above it; the data
# rows in this sheet carry no explicit cell style, so strip that back off.
$ws.Rows.Item(2).ClearFormats()

# New top row: one day after the previous newest date (2025-12-30 -> 2025-12-31).
# Force text so Excel doesn't auto-convert the date-shaped string to a serial
# date number, then drop the helper format so the cell stays style-less,
# matching every other date cell in the column.
$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "2025-12-31"
$ws.Cells.Item(2, 1).ClearFormats()
$ws.Cells.Item(2, 2).Value = 783.5
$ws.Cells.Item(2, 3).Value = 1112
$ws.Cells.Item(2, 4).Value = 3610
